$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("Value of Approved Payments" under Direct Debit Payments) currently
# mirrors the blank "Dishonored Direct Debits" style (s=5). Re-style A12 to
# match the other "Value/Count of Approved Payments" rows (e.g. A7, s=4) by
# copying the formatting across, then write the calculated approved Direct
# Debit payment values for columns B through N.

$ws.Range("A7").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B12:N12").Style = "Normal"

$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 92390.01700000002
$ws.Range("G12").Value = 41463.556000000004
$ws.Range("H12").Value = 86566.3205
$ws.Range("I12").Value = 29267.3985
$ws.Range("J12").Value = 62031.90700000001
$ws.Range("K12").Value = 49470.471
$ws.Range("L12").Value = 36927.276
$ws.Range("M12").Value = 64959.96200000001
$ws.Range("N12").Value = 60838.859500000006
